$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Residence address: Hong Kong address -> US address
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Flat C 16th Floor Tower 10, Park Central, 9 Tong Tak Street, Tseung Kwan O",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1082 Mt Dana Dr, Chula Vista, CA 91913", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Mobile number: HK number -> US number (with trailing PDF mark U+202C)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "+852.9732.6715",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "+1-619-800-0859" + [char]0x202C, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Remove the stray "_GoBack" bookmark left over from editing
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 4. Job title: "Data Science Manager, ..." -> "Senior Manager - Data
#    Science, ..." (all bold, same size/color as surrounding text)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Data Science Manager, Global Supply Chain Analytics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Senior Manager " + [char]0x2013 + " Data Science, Global Supply Chain Analytics", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Company line: "LF Logistics - A Li & Fung Company (Hong Kong)"
#    -> "LF Logistics (Global)"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("LF Logistics " + [char]0x2013 + " A Li & Fung Company (Hong Kong)")
if ($found) {
    $base = $rng.Start
    # Edit back-to-front so earlier offsets stay valid.
    $r3 = $d.Range($base + 34, $base + 46)
    $r3.Text = ")"
    $r3.Font.Color = 0
    $r3.Font.Bold = $true

    $r2 = $d.Range($base + 14, $base + 34)
    $r2.Text = "Global"
    $r2.Font.Color = 0
    $r2.Font.Bold = $true

    $r1 = $d.Range($base + 13, $base + 14)
    $r1.Text = "("
    $r1.Font.Color = 0
    $r1.Font.Bold = $true
}

# ---------------------------------------------------------------------
# 6. "Responsible for optimizing Li & Fung's" -> "...LF's"
#    "Developed the data science platform for Li & Fung's" -> "...LF's"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Li & Fung" + [char]0x2019 + "s",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LF" + [char]0x2019 + "s", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. "... on premise or in the cloud ..." gramStart/gramEnd markers
#    around "or" do not change the visible text, nothing to do.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 8. "IPG Mediabrands (Philippines)" -> "IPG Mediabrands (Asia Pacific)"
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Mediabrands (Philippines)")
if ($found2) {
    $base2 = $rng2.Start
    $sub = $d.Range($base2 + 11, $base2 + 25)
    $sub.Text = " (Asia Pacific)"
}
